$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title: "sprint2" -> "sprint4"
# ------------------------------------------------------------------
$d.Content.Find.Execute("reunião da retrospectiva sprint2", $true, $false, $false, $false, $false, `
    $true, 1, $false, "reunião da retrospectiva sprint4", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Meeting date "REUNIÃO (04" -> "REUNIÃO (11", and move the
#    "_GoBack" bookmark so it sits right after that run (it used to
#    live at the very end of the document).
# ------------------------------------------------------------------
$d.Content.Find.Execute("REUNIÃO (04", $true, $false, $false, $false, $false, `
    $true, 1, $false, "REUNIÃO (11", 2) | Out-Null

$found = $d.Content
$found.Find.Execute("REUNIÃO (11") | Out-Null
$bmRange = $d.Range($found.End, $found.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 3) "O que foi bem" paragraph body text
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "equipe teve uma boa comunicação, sendo assim foi definido a divisão nessa Sprint por história, por essa razão cada membro da equipe desenvolveu mais responsabilidade pelas tarefas.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "equipe conseguiu seguir a carga horária diária de acordo com as horas propostas pela Sprint BurnDown.", 2) | Out-Null

# ------------------------------------------------------------------
# 4) "O que não foi tão bem" paragraph body text
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "A equipe teve um pouco de dificuldade em dividir a carga horária diariamente para trabalhar com o projeto, devido a inúmeras tarefas que foram definidas na Sprint-3.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "A equipe teve dificuldade com o controle de versões do sistema, onde ficou impossibilitado o compartilhamento do código entre a equipe.", 2) | Out-Null

# ------------------------------------------------------------------
# 5) "Melhoras" paragraph: drop the hanging 1440-twip left indent and
#    rewrite both runs of body text.
# ------------------------------------------------------------------
$headPara = $d.Content
$headPara.Find.Execute("Melhoras") | Out-Null
$bodyPara = $headPara.Paragraphs(1).Next()
$bodyPara.LeftIndent = 0

$d.Content.Find.Execute("A equipe busca ", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "                    A equipe busca trabalhar com tarefas menores que possibilitem o compartilhamento do código rapidamente, com isso evitar erros de mesclagem ", 2) | Out-Null

$d.Content.Find.Execute( `
    "uma melhor performance com o modelo Scrum, sem dificuldades e um desenvolvimento rápido e eficiente.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "no controle de versões.", 2) | Out-Null
